$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New accelerometer readings: rows 2-31 (A:C). Existing A2:C21 data is replaced
# and the sheet grows to A1:C31 (10 additional rows of data appended on May 9th).
$data = New-Object 'object[,]' 30,3
$data[0,0] = -12.57026290893555
$data[0,1] = -8.004177093505859
$data[0,2] = -9.469242095947266
$data[1,0] = -11.83102798461914
$data[1,1] = -4.641905784606934
$data[1,2] = -9.06204319000244
$data[2,0] = 5.961018562316895
$data[2,1] = -11.27237701416016
$data[2,2] = 6.515813827514648
$data[3,0] = 5.790210247039795
$data[3,1] = -16.61598014831543
$data[3,2] = 5.962498664855957
$data[4,0] = -44.78383636474609
$data[4,1] = -8.077349662780762
$data[4,2] = -27.36338424682617
$data[5,0] = -42.66062545776367
$data[5,1] = -5.792407035827637
$data[5,2] = -30.21708297729492
$data[6,0] = -9.564473152160645
$data[6,1] = -21.24031639099121
$data[6,2] = -0.3582277297973633
$data[7,0] = -10.97329330444336
$data[7,1] = -18.38039779663086
$data[7,2] = 0.1557941436767578
$data[8,0] = 9.481355667114258
$data[8,1] = -9.28108024597168
$data[8,2] = 9.836421966552734
$data[9,0] = 10.69957828521728
$data[9,1] = -16.89757537841797
$data[9,2] = 15.06509208679199
$data[10,0] = -12.81566715240478
$data[10,1] = -8.839512825012207
$data[10,2] = -14.14643669128418
$data[11,0] = -8.493080139160156
$data[11,1] = -3.757828950881958
$data[11,2] = -17.9176139831543
$data[12,0] = -10.02292919158936
$data[12,1] = -2.351483345031738
$data[12,2] = 5.807761192321777
$data[13,0] = -12.75368976593018
$data[13,1] = -6.235836982727051
$data[13,2] = 6.349725723266602
$data[14,0] = -79.57471466064453
$data[14,1] = -37.79908752441406
$data[14,2] = -28.62195587158203
$data[15,0] = -83.42860412597656
$data[15,1] = -34.91292572021484
$data[15,2] = -30.07511520385743
$data[16,0] = -3.41340184211731
$data[16,1] = -14.34111213684082
$data[16,2] = -5.178817749023438
$data[17,0] = -5.549047946929932
$data[17,1] = -10.40872192382812
$data[17,2] = -4.603366851806641
$data[18,0] = 13.1239423751831
$data[18,1] = -7.421818733215332
$data[18,2] = 9.609323501586914
$data[19,0] = 13.5955753326416
$data[19,1] = -13.69983100891113
$data[19,2] = 12.63335609436035
$data[20,0] = 3.66697096824646
$data[20,1] = -4.846949577331543
$data[20,2] = -18.3593635559082
$data[21,0] = 8.207674980163574
$data[21,1] = -0.2672674059867859
$data[21,2] = -20.82013320922852
$data[22,0] = -3.444841146469116
$data[22,1] = 0.0402781963348388
$data[22,2] = -6.625090599060059
$data[23,0] = -3.938364028930664
$data[23,1] = -3.24791407585144
$data[23,2] = -6.832444190979004
$data[24,0] = -71.62106323242188
$data[24,1] = -3.320010900497437
$data[24,2] = -19.81892585754395
$data[25,0] = -75.16300201416016
$data[25,1] = -2.210659265518188
$data[25,2] = -18.6042366027832
$data[26,0] = -11.13384056091309
$data[26,1] = -14.18883228302002
$data[26,2] = -7.074286460876465
$data[27,0] = -9.418439865112305
$data[27,1] = -9.676417350769045
$data[27,2] = -8.096508026123047
$data[28,0] = 17.99904441833496
$data[28,1] = 10.46211528778076
$data[28,2] = 5.739192962646484
$data[29,0] = 21.18538665771484
$data[29,1] = 4.166534900665283
$data[29,2] = 8.056502342224121

$ws.Range("A2:C31").Value = $data
